$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021年" row (row 5) to the "按登记注册类型分..." table,
# matching the style used for the existing year-label cells in column A.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "2021年"

$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 162708
$ws.Range("D5").Value = 861
$ws.Range("E5").Value = 1665
$ws.Range("F5").Value = 8463
$ws.Range("G5").Value = 25972
$ws.Range("H5").Value = 7382
$ws.Range("I5").Value = 130322
$ws.Range("J5").Value = 36
$ws.Range("K5").Value = 230
$ws.Range("L5").Value = 5054
$ws.Range("M5").Value = 208
